# journal de travail update
# Fills in row 16 (12.03.2020, 5h, new comment, "1h30") and row 17's date
# (13.03.2020), matching the existing look-and-feel of the rows above
# (same date style, same wrap-text comment style), then moves the
# selection to E15 as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 -----------------------------------------------------------
# A16: date 12.03.2020 (serial 43902) - copy the date format already used
# on A13 (numFmtId 14, centered, no border) then set the value.
$ws.Range("A13").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 43902

# B16: week number
$ws.Range("B16").Value = 5

# C16: task comment - copy the wrap-text style already used on E16 so the
# text wraps inside the cell, then set the text itself (a new, unique
# string).
$ws.Range("E16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Affichage des bateaux coulés et messages de victoire, le nombre de cases ratés quand tous les bateaux sont coulés"

# D16: duration, reuses the already-existing "1h30" text used elsewhere
# in the sheet.
$ws.Range("D16").Value = "1h30"

# Row height grows to fit the new wrapped comment.
$ws.Rows.Item(16).RowHeight = 45

# --- Row 17 -------------------------------------------------------------
# A17: date 13.03.2020 (serial 43903), same date style as above.
$ws.Range("A13").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 43903

# --- View state ---------------------------------------------------------
$ws.Range("E15").Select()
